$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.508.89"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.238.43"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "'306.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'94.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("D7").Value = "'0.570"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "'34.82"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.28%  "
$ws.Range("D11").Value = "'0.0802"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("D12").Value = "'7.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "2.272.54"
$ws.Range("E14").Value = "  +1.65%  "
$ws.Range("D15").Value = "'0.832"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.14%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "44.216.32"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "'6.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "'11.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "'65.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "'236.67"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("D26").Value = "'2.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.03%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").Value = "'9.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.13%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'37.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").Value = "'5.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("D30").Value = "'19.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.08%  "
$ws.Range("D31").Value = "'152.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").Value = "'2.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("B34").Value = "LidoDAOToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D34").Value = "'3.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.63%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.110"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("E37").Value = "  +2.87%  "
$ws.Range("D38").Value = "'15.13"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.95%  "
$ws.Range("D39").Value = "'3.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").Value = "'3.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").Value = "'0.0299"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "1.798.22"
$ws.Range("E43").Value = "  +4.38%  "
$ws.Range("D44").Value = "'1.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.18%  "
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("D46").Value = "'78.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.70%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "'70.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.38%  "
$ws.Range("D48").Value = "'98.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("B49").Value = "THORChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D49").Value = "'4.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.92%  "
$ws.Range("D50").Value = "'8.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("D51").Value = "'54.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
